# Auto-generated Excel COM-interop edit script
# Applies scheduled-runner price/profit updates to the Famfrit_Profits workbook
# (columns H:N across ALC, ARM, CRP, CUL, GSM, LTW, WVR sheets)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4: Root Rush
$ws.Range("H4").Value = 332.5
$ws.Range("I4").Value = 332.5
$ws.Range("K4").Value = 332.5
$ws.Range("M4").Value = -218.5

# Row 9: Distill, My Heart
$ws.Range("H9").Value = 7548.6875
$ws.Range("I9").Value = 14463.429
$ws.Range("K9").Value = 14463.429
$ws.Range("M9").Value = -14294.429

# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 10995789
$ws.Range("I62").Value = 5280.25
$ws.Range("K62").Value = 5280.25
$ws.Range("M62").Value = -4656.25

# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 10995789
$ws.Range("I65").Value = 5280.25
$ws.Range("K65").Value = 26401.25
$ws.Range("M65").Value = -23281.25

# Row 69: Steeling the Knife, Steeling the Mind
$ws.Range("H69").Value = 4600
$ws.Range("I69").Value = 4600
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 13800
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -12926
$ws.Range("N69").Value = $null

# Row 72: Surgical Substitution (L)
$ws.Range("H72").Value = 4600
$ws.Range("I72").Value = 4600
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 41400
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -37032
$ws.Range("N72").Value = $null

# Row 76: Warding Off Temptation
$ws.Range("H76").Value = 7591.2144
$ws.Range("I76").Value = 7253.857
$ws.Range("K76").Value = 7253.857
$ws.Range("M76").Value = -6938.857

# Row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 7591.2144
$ws.Range("I79").Value = 7253.857
$ws.Range("K79").Value = 7253.857
$ws.Range("M79").Value = -6161.857

# Row 99: Rumor Has It
$ws.Range("H99").Value = 1016.1429
$ws.Range("I99").Value = 400
$ws.Range("J99").Value = 1118.8334
$ws.Range("K99").Value = 1200
$ws.Range("L99").Value = 3356.5002
$ws.Range("M99").Value = 298
$ws.Range("N99").Value = -6352.5002

# Row 113: Amaro Kart
$ws.Range("H113").Value = 7789.6
$ws.Range("I113").Value = 7800
$ws.Range("J113").Value = 7748
$ws.Range("K113").Value = 7800
$ws.Range("L113").Value = 7748
$ws.Range("M113").Value = -4546
$ws.Range("N113").Value = -14256

# Row 138: All-night Crafting
$ws.Range("H138").Value = 12106.134
$ws.Range("I138").Value = 3198.4
$ws.Range("J138").Value = 16560
$ws.Range("K138").Value = 9595.200000000001
$ws.Range("L138").Value = 49680
$ws.Range("M138").Value = -4455.200000000001
$ws.Range("N138").Value = -59960

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 1366.7262
$ws.Range("I32").Value = 936.452
$ws.Range("K32").Value = 936.452
$ws.Range("M32").Value = -649.452

# Row 37: Get Shirty
$ws.Range("H37").Value = 20898.625
$ws.Range("I37").Value = 17999.834
$ws.Range("J37").Value = 29595
$ws.Range("K37").Value = 17999.834
$ws.Range("L37").Value = 29595
$ws.Range("M37").Value = -17726.834
$ws.Range("N37").Value = -30141

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 23257896
$ws.Range("I74").Value = 33334900
$ws.Range("K74").Value = 33334900
$ws.Range("M74").Value = -33334026

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 23257896
$ws.Range("I77").Value = 33334900
$ws.Range("K77").Value = 166674500
$ws.Range("M77").Value = -166670132

# Row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 3724.75
$ws.Range("I102").Value = 1966.3334
$ws.Range("J102").Value = 9000
$ws.Range("K102").Value = 1966.3334
$ws.Range("L102").Value = 9000
$ws.Range("M102").Value = -344.3334
$ws.Range("N102").Value = -12244

# Row 113: Catching an Earful
$ws.Range("H113").Value = 398
$ws.Range("J113").Value = 398
$ws.Range("L113").Value = 398
$ws.Range("N113").Value = -9076

$ws = $wb.Worksheets.Item("CRP")
# Row 10: Spears and Sorcery
$ws.Range("H10").Value = 964.5
$ws.Range("I10").Value = 224.9
$ws.Range("J10").Value = 1889
$ws.Range("K10").Value = 224.9
$ws.Range("L10").Value = 1889
$ws.Range("M10").Value = -85.90000000000001
$ws.Range("N10").Value = -2167

# Row 31: Wall Not Found
$ws.Range("H31").Value = 7105.8887
$ws.Range("I31").Value = 5790.8
$ws.Range("J31").Value = 8749.75
$ws.Range("K31").Value = 5790.8
$ws.Range("L31").Value = 8749.75
$ws.Range("M31").Value = -5495.8
$ws.Range("N31").Value = -9339.75

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 7105.8887
$ws.Range("I34").Value = 5790.8
$ws.Range("J34").Value = 8749.75
$ws.Range("K34").Value = 5790.8
$ws.Range("L34").Value = 8749.75
$ws.Range("M34").Value = -5588.8
$ws.Range("N34").Value = -9153.75

# Row 99: O Pine
$ws.Range("H99").Value = 8165.6
$ws.Range("I99").Value = 2725
$ws.Range("J99").Value = 9201.904
$ws.Range("K99").Value = 2725
$ws.Range("L99").Value = 9201.904
$ws.Range("M99").Value = -1227
$ws.Range("N99").Value = -12197.904

# Row 126: A Better Conductor
$ws.Range("H126").Value = 8165.6
$ws.Range("I126").Value = 2725
$ws.Range("J126").Value = 9201.904
$ws.Range("K126").Value = 8175
$ws.Range("L126").Value = 27605.712
$ws.Range("M126").Value = -5705
$ws.Range("N126").Value = -32545.712

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 2996.4583
$ws.Range("I132").Value = 2677.9546
$ws.Range("K132").Value = 8033.8638
$ws.Range("M132").Value = -5503.8638

$ws = $wb.Worksheets.Item("CUL")
# Row 34: Fever Pitch
$ws.Range("H34").Value = 379
$ws.Range("I34").Value = 99.5
$ws.Range("J34").Value = 565.3333
$ws.Range("K34").Value = 298.5
$ws.Range("L34").Value = 1695.9999
$ws.Range("M34").Value = -214.5
$ws.Range("N34").Value = -1863.9999

# Row 39: Bloody Good Tart, This
$ws.Range("H39").Value = 3130.4546
$ws.Range("I39").Value = 1491.2858
$ws.Range("J39").Value = 5999
$ws.Range("K39").Value = 4473.857400000001
$ws.Range("L39").Value = 17997
$ws.Range("M39").Value = -4179.857400000001
$ws.Range("N39").Value = -18585

# Row 55: Pagan Pastries
$ws.Range("H55").Value = 1003266.7
$ws.Range("I55").Value = 1801640
$ws.Range("K55").Value = 5404920
$ws.Range("M55").Value = -5404743

# Row 88: Don't Let It Fall Apart
$ws.Range("H88").Value = 3000
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = $null

# Row 91: Better Come Back with a Sandwich (L)
$ws.Range("H91").Value = 3000
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = $null

$ws = $wb.Worksheets.Item("GSM")
# Row 7: Water of Life
$ws.Range("H7").Value = 76000000
$ws.Range("I7").Value = 76000000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 76000000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -75999888
$ws.Range("N7").Value = $null

# Row 8: Gods of Small Things
$ws.Range("H8").Value = 76000000
$ws.Range("I8").Value = 76000000
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 76000000
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -75999861
$ws.Range("N8").Value = $null

# Row 11: A Ringing Success
$ws.Range("H11").Value = 12022615
$ws.Range("I11").Value = 22302570
$ws.Range("J11").Value = 29333.334
$ws.Range("K11").Value = 22302570
$ws.Range("L11").Value = 29333.334
$ws.Range("M11").Value = -22302431
$ws.Range("N11").Value = -29611.334

# Row 13: A Needle Is a Small Sword
$ws.Range("H13").Value = 220
$ws.Range("I13").Value = 136
$ws.Range("J13").Value = 325
$ws.Range("K13").Value = 136
$ws.Range("L13").Value = 325
$ws.Range("M13").Value = 3
$ws.Range("N13").Value = -603

# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 9837.333000000001
$ws.Range("I70").Value = 7256
$ws.Range("J70").Value = 15000
$ws.Range("K70").Value = 7256
$ws.Range("L70").Value = 15000
$ws.Range("M70").Value = -6986
$ws.Range("N70").Value = -15540

# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 9837.333000000001
$ws.Range("I73").Value = 7256
$ws.Range("J73").Value = 15000
$ws.Range("K73").Value = 7256
$ws.Range("L73").Value = 15000
$ws.Range("M73").Value = -6320
$ws.Range("N73").Value = -16872

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 1970.5385
$ws.Range("I122").Value = 1662.5
$ws.Range("J122").Value = 2997.3333
$ws.Range("K122").Value = 4987.5
$ws.Range("L122").Value = 8991.999899999999
$ws.Range("M122").Value = -2537.5
$ws.Range("N122").Value = -13891.9999

# Row 132: On Board for Lar
$ws.Range("H132").Value = 5971.75
$ws.Range("I132").Value = 5307.3716
$ws.Range("J132").Value = 8555.444
$ws.Range("K132").Value = 15922.1148
$ws.Range("L132").Value = 25666.332
$ws.Range("M132").Value = -13392.1148
$ws.Range("N132").Value = -30726.332

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 2141.5
$ws.Range("I22").Value = 1637.375
$ws.Range("J22").Value = 3149.75
$ws.Range("K22").Value = 1637.375
$ws.Range("L22").Value = 3149.75
$ws.Range("M22").Value = -1342.375
$ws.Range("N22").Value = -3739.75

# Row 27: Fire and Hide
$ws.Range("H27").Value = 2141.5
$ws.Range("I27").Value = 1637.375
$ws.Range("J27").Value = 3149.75
$ws.Range("K27").Value = 1637.375
$ws.Range("L27").Value = 3149.75
$ws.Range("M27").Value = -1530.375
$ws.Range("N27").Value = -3363.75

# Row 40: Best Served Toad
$ws.Range("H40").Value = 6348.1
$ws.Range("I40").Value = 6211.2104
$ws.Range("J40").Value = 8949
$ws.Range("K40").Value = 6211.2104
$ws.Range("L40").Value = 8949
$ws.Range("M40").Value = -6075.2104
$ws.Range("N40").Value = -9221

# Row 82: Trainin' the Neck
$ws.Range("H82").Value = 2608.8
$ws.Range("I82").Value = 1977.3
$ws.Range("J82").Value = 3871.8
$ws.Range("K82").Value = 1977.3
$ws.Range("L82").Value = 3871.8
$ws.Range("M82").Value = -1616.3
$ws.Range("N82").Value = -4593.8

# Row 85: Training Is Only Skintight (L)
$ws.Range("H85").Value = 2608.8
$ws.Range("I85").Value = 1977.3
$ws.Range("J85").Value = 3871.8
$ws.Range("K85").Value = 1977.3
$ws.Range("L85").Value = 3871.8
$ws.Range("M85").Value = -729.3
$ws.Range("N85").Value = -6367.8

# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 3237.7144
$ws.Range("I93").Value = 2871.5
$ws.Range("J93").Value = 7998.5
$ws.Range("K93").Value = 2871.5
$ws.Range("L93").Value = 7998.5
$ws.Range("M93").Value = -1623.5
$ws.Range("N93").Value = -10494.5

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins
$ws.Range("H132").Value = 3302.7585
$ws.Range("I132").Value = 2865.652
$ws.Range("K132").Value = 8596.956
$ws.Range("M132").Value = -6066.956
